# DaySale_2025-06-18_00-00.xlsx update:
#  - a new product row ("face mask" / "ماسك وجه") is inserted just above the
#    totals row at the bottom of the sheet
#  - the totals row's grand total is bumped up by the new line's sale total
#  - the footer row (timestamp / page no. / developer credit) shifts down by
#    one row and its timestamp is refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow   = 42
$sumRow   = 43
$footRow  = 44

# ---------------------------------------------------------------------------
# 1) Push the totals row (and everything below it) down by one row so there
#    is a blank row ready to receive the new product line.
# ---------------------------------------------------------------------------
$ws.Rows.Item($newRow).Insert()

# ---------------------------------------------------------------------------
# 2) Re-create the merged cell groups for the new data row (same layout as
#    every other product row: A:B, C:G, H:K, L:M, N:O, P and Q stand alone).
# ---------------------------------------------------------------------------
$ws.Range("A" + $newRow + ":B" + $newRow).Merge()
$ws.Range("C" + $newRow + ":G" + $newRow).Merge()
$ws.Range("H" + $newRow + ":K" + $newRow).Merge()
$ws.Range("L" + $newRow + ":M" + $newRow).Merge()
$ws.Range("N" + $newRow + ":O" + $newRow).Merge()

# ---------------------------------------------------------------------------
# 3) Clone the formatting of the row above (row 41) onto every column group
#    of the new row so fonts/fills/borders/number-formats match the rest of
#    the table exactly.
# ---------------------------------------------------------------------------
$prevRow = $newRow - 1

$ws.Range("A" + $prevRow + ":B" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":B" + $newRow).PasteSpecial(-4122)

$ws.Range("C" + $prevRow + ":G" + $prevRow).Copy()
$ws.Range("C" + $newRow + ":G" + $newRow).PasteSpecial(-4122)

$ws.Range("H" + $prevRow + ":K" + $prevRow).Copy()
$ws.Range("H" + $newRow + ":K" + $newRow).PasteSpecial(-4122)

$ws.Range("L" + $prevRow + ":M" + $prevRow).Copy()
$ws.Range("L" + $newRow + ":M" + $newRow).PasteSpecial(-4122)

$ws.Range("N" + $prevRow + ":O" + $prevRow).Copy()
$ws.Range("N" + $newRow + ":O" + $newRow).PasteSpecial(-4122)

$ws.Range("P" + $prevRow).Copy()
$ws.Range("P" + $newRow).PasteSpecial(-4122)

$ws.Range("Q" + $prevRow).Copy()
$ws.Range("Q" + $newRow).PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Row height matches the other "tall" data rows.
$ws.Rows.Item($newRow).RowHeight = 25.5

# ---------------------------------------------------------------------------
# 4) Fill in the new product's data.
# ---------------------------------------------------------------------------
$ws.Range("A" + $newRow).Value = 36          # م (line number)
$ws.Range("C" + $newRow).Value = "ماسك وجه"   # الاسم (name)
$ws.Range("H" + $newRow).Value = "13:0"       # الرصيد الحالي (current balance)

# L (order limit) re-uses the same literal "0" already used throughout the
# sheet - copy it straight from L41 so it stays a shared text value, not a
# recalculated number.
$ws.Range("L" + $prevRow).Copy()
$ws.Range("L" + $newRow).PasteSpecial(-4163)

$ws.Range("N" + $newRow).Value = "30.00"      # السعر (price)

# P (sale price) needs the literal text "30.0000" preserved exactly (not
# reduced to the number 30). Stage it in a throw-away column that is plain
# text-formatted, then bring only the *value* across so the real cell's
# number format (inherited above) is left untouched.
$stageCol = 100
$ws.Cells.Item(1, $stageCol).NumberFormat = "@"
$ws.Cells.Item(1, $stageCol).Value = "30.0000"
$ws.Cells.Item(1, $stageCol).Copy()
$ws.Range("P" + $newRow).PasteSpecial(-4163)
$ws.Columns.Item($stageCol).Delete()

# Q (transaction count) re-uses the same "1:0" literal already present on
# row 7 - copy it so the value lands as the identical shared string.
$ws.Range("Q7").Copy()
$ws.Range("Q" + $newRow).PasteSpecial(-4163)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Totals row: bump the grand total by the new line's sale total (30) and
#    match the auto-fit height Excel produced once the table grew by a row.
# ---------------------------------------------------------------------------
$ws.Range("P" + $sumRow).Value = $ws.Range("P" + $sumRow).Value + 30
$ws.Rows.Item($sumRow).RowHeight = 24.75

# ---------------------------------------------------------------------------
# 6) Footer row: refresh the printed timestamp to the new save time.
# ---------------------------------------------------------------------------
$ws.Range("A" + $footRow).Value = "Wednesday, 18 June, 2025 1:29 PM"
